$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 41933.2960456292
$ws.Cells.Item(2, 5).Value = 98156.4036300151
$ws.Cells.Item(2, 6).Value = 117714.991337292
$ws.Cells.Item(2, 9).Value = 13578.2960456292

$ws.Cells.Item(3, 2).Value = 28498.9131744028
$ws.Cells.Item(3, 5).Value = 116469.059241186
$ws.Cells.Item(3, 6).Value = 147374.878940713
$ws.Cells.Item(3, 9).Value = 25235.9131744028

$ws.Cells.Item(4, 2).Value = 49261.2437854625
$ws.Cells.Item(4, 5).Value = 123335.073555769
$ws.Cells.Item(4, 6).Value = 152732.312531613
$ws.Cells.Item(4, 9).Value = 46392.2437854625

$ws.Cells.Item(5, 2).Value = 130638.506378801
$ws.Cells.Item(5, 9).Value = 127637.506378801

$ws.Cells.Item(6, 2).Value = 253992.504175672
$ws.Cells.Item(6, 9).Value = 250073.504175672

$ws.Cells.Item(7, 2).Value = 297047.660429333
$ws.Cells.Item(7, 9).Value = 290787.660429333

$ws.Cells.Item(8, 2).Value = 257853.623716744
$ws.Cells.Item(8, 9).Value = 198710.623716744

$ws.Cells.Item(9, 2).Value = 165149.199836224
$ws.Cells.Item(9, 5).Value = 352465.143848049
$ws.Cells.Item(9, 9).Value = 106006.199836224

$ws.Cells.Item(10, 2).Value = 156186.493974779
$ws.Cells.Item(10, 9).Value = 41016.4939747785

$ws.Cells.Item(11, 2).Value = 139594.218504904
$ws.Cells.Item(11, 9).Value = -60280.7814950955

$ws.Cells.Item(12, 2).Value = 116923.176838118
$ws.Cells.Item(12, 9).Value = -72093.8231618823

$ws.Cells.Item(13, 2).Value = 87422.0541342418
$ws.Cells.Item(13, 9).Value = -44375.9458657582

$ws.Cells.Item(14, 2).Value = 52286.2157739597
$ws.Cells.Item(14, 9).Value = -10922.7842260403

$ws.Cells.Item(15, 2).Value = 31356.6137338612
$ws.Cells.Item(15, 9).Value = 12827.6137338612

$ws.Cells.Item(16, 2).Value = 60784.7766028931
$ws.Cells.Item(16, 9).Value = 24578.7766028931

$ws.Cells.Item(17, 2).Value = 153742.057868615
$ws.Cells.Item(17, 9).Value = 9575.05786861462

$ws.Cells.Item(18, 2).Value = 284405.41632631
$ws.Cells.Item(18, 9).Value = 43567.4163263101

$ws.Cells.Item(19, 2).Value = 323520.857783823
$ws.Cells.Item(19, 9).Value = 99094.8577838232

$ws.Cells.Item(20, 2).Value = 274630.806402726
$ws.Cells.Item(20, 9).Value = 111643.806402726

$ws.Cells.Item(21, 2).Value = 172152.672434554
$ws.Cells.Item(21, 9).Value = 100869.672434554

$ws.Cells.Item(22, 2).Value = 170040.510184414
$ws.Cells.Item(22, 9).Value = 83405.5101844135

$ws.Cells.Item(23, 2).Value = 152333.747875243
$ws.Cells.Item(23, 9).Value = 20263.7478752426

$ws.Cells.Item(24, 2).Value = 130298.722189757
$ws.Cells.Item(24, 9).Value = 16825.7221897574

$ws.Cells.Item(25, 2).Value = 98304.4437248791
$ws.Cells.Item(25, 9).Value = 23883.4437248791

$ws.Cells.Item(26, 2).Value = 59480.2635064802
$ws.Cells.Item(26, 9).Value = 15297.2635064802

$ws.Cells.Item(27, 2).Value = 35673.4705155926
$ws.Cells.Item(27, 9).Value = 17539.4705155926

$ws.Cells.Item(28, 2).Value = 70780.9692509437
$ws.Cells.Item(28, 9).Value = 37493.9692509437

$ws.Cells.Item(29, 2).Value = 172271.195784828
$ws.Cells.Item(29, 9).Value = 127477.195784828

$ws.Cells.Item(30, 2).Value = 309625.428674793
$ws.Cells.Item(30, 9).Value = 239583.428674793

$ws.Cells.Item(31, 2).Value = 347869.669583153
$ws.Cells.Item(31, 9).Value = 207208.669583153

$ws.Cells.Item(32, 2).Value = 293017.113851282
$ws.Cells.Item(32, 9).Value = 184044.113851282

$ws.Cells.Item(33, 2).Value = 185403.566830394
$ws.Cells.Item(33, 9).Value = 130875.566830394

$ws.Cells.Item(34, 2).Value = 183652.463428508
$ws.Cells.Item(34, 9).Value = 138780.463428508

$ws.Cells.Item(35, 2).Value = 167378.034410658
$ws.Cells.Item(35, 9).Value = 119983.034410658

$ws.Cells.Item(36, 2).Value = 146040.031292512
$ws.Cells.Item(36, 9).Value = 95407.0312925121

$ws.Cells.Item(37, 2).Value = 112497.089271181
$ws.Cells.Item(37, 9).Value = 84750.0892711806

$ws.Cells.Item(38, 2).Value = 69773.5028059606
$ws.Cells.Item(38, 9).Value = 67289.5028059606

$ws.Cells.Item(39, 2).Value = 42654.730055987
$ws.Cells.Item(39, 9).Value = 38719.730055987

$ws.Cells.Item(40, 2).Value = 81630.2068740226
$ws.Cells.Item(40, 9).Value = 71999.2068740226

$ws.Cells.Item(41, 2).Value = 190275.201586783
$ws.Cells.Item(41, 9).Value = 170170.201586783

$ws.Cells.Item(42, 2).Value = 333774.474146971
$ws.Cells.Item(42, 9).Value = 242515.474146971

$ws.Cells.Item(43, 2).Value = 373058.665848598
$ws.Cells.Item(43, 9).Value = -60025.3341514025

$ws.Cells.Item(44, 2).Value = 315788.636438919
$ws.Cells.Item(44, 9).Value = -141423.363561081

$ws.Cells.Item(45, 2).Value = 202069.475104455
$ws.Cells.Item(45, 9).Value = 8531.47510445534

$ws.Cells.Item(46, 2).Value = 200075.090750977
$ws.Cells.Item(46, 9).Value = 33095.0907509766

$ws.Cells.Item(47, 2).Value = 184280.379858738
$ws.Cells.Item(47, 9).Value = 18753.3798587384

$ws.Cells.Item(48, 2).Value = 163225.843194668
$ws.Cells.Item(48, 9).Value = 70270.8431946682

$ws.Cells.Item(49, 2).Value = 128485.497160848
$ws.Cells.Item(49, 9).Value = 82335.4971608482
